# Attempt to fix the micron detection issue.
# Recomputed derived measurements (length/width/area/etc.) for the first
# three particles after correcting the pixel-size / micron detection,
# and tightened a couple of column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (particle id=1): H..R recomputed ---
$ws.Cells.Item(2,8).Value2  = 192210
$ws.Cells.Item(2,9).Value2  = 400.56785152729407
$ws.Cells.Item(2,10).Value2 = 126020.75103673227
$ws.Cells.Item(2,11).Value2 = 213.80144697404555
$ws.Cells.Item(2,12).Value2 = 5335.2226720647768
$ws.Cells.Item(2,13).Value2 = 0.0556348379808402
$ws.Cells.Item(2,14).Value2 = 0.078128335136030333
$ws.Cells.Item(2,15).Value2 = 334.58405389938088
$ws.Cells.Item(2,16).Value2 = 549.73586702044634
$ws.Cells.Item(2,17).Value2 = 26.565763847707888
$ws.Cells.Item(2,18).Value2 = 26.565763847707888

# --- Row 3 (particle id=2): E..R recomputed ---
$ws.Cells.Item(3,5).Value2  = 295.54655870445345
$ws.Cells.Item(3,6).Value2  = 225.91093117408906
$ws.Cells.Item(3,7).Value2  = 1.3082437275985663
$ws.Cells.Item(3,8).Value2  = 49529
$ws.Cells.Item(3,9).Value2  = 203.33768814013857
$ws.Cells.Item(3,10).Value2 = 32473.241652870885
$ws.Cells.Item(3,11).Value2 = 82.321156136135841
$ws.Cells.Item(3,12).Value2 = 1429.9595141700404
$ws.Cells.Item(3,13).Value2 = 0.19956667663600342
$ws.Cells.Item(3,14).Value2 = 0.34628551744154062
$ws.Cells.Item(3,15).Value2 = 839.44668779906726
$ws.Cells.Item(3,16).Value2 = 468.32403238506731
$ws.Cells.Item(3,17).Value2 = 30.050364680042847
$ws.Cells.Item(3,18).Value2 = 30.050364680042847

# --- Row 4 (particle id=3): E..R recomputed ---
$ws.Cells.Item(4,5).Value2  = 9.7165991902834001
$ws.Cells.Item(4,6).Value2  = 6.4777327935222671
$ws.Cells.Item(4,7).Value2  = 1.5
$ws.Cells.Item(4,8).Value2  = 110
$ws.Cells.Item(4,9).Value2  = 9.5826239234059862
$ws.Cells.Item(4,10).Value2 = 72.120506810470573
$ws.Cells.Item(4,11).Value2 = 3.5330940923602463
$ws.Cells.Item(4,12).Value2 = 31.247639500791408
$ws.Cells.Item(4,13).Value2 = 0.92818426720588398
$ws.Cells.Item(4,14).Value2 = 0.18568439238296125
$ws.Cells.Item(4,15).Value2 = 317.92727272727274
$ws.Cells.Item(4,16).Value2 = 616.25454545454545
$ws.Cells.Item(4,17).Value2 = 6.851449106449107
$ws.Cells.Item(4,18).Value2 = 6.851449106449107

# --- Column width tweaks ---
# Column B: 3.140625 -> 2.85546875 (closest reachable pixel-snapped width)
$ws.Columns.Item(2).ColumnWidth = 2.0
# Column M: 12.7109375 -> 13.7109375 (closest reachable pixel-snapped width)
$ws.Columns.Item(13).ColumnWidth = 12.833333333333332
